$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The gage-trigger list is being trimmed: three USACE gages (82770,
# 76220, 76593) are no longer triggers and their rows are removed
# entirely, with every row below shifting up to close the gaps.
$gagesToRemove = @("82770", "76220", "76593")

# Find the current row number for each gage (column B) and delete them
# bottom-to-top so earlier row numbers stay valid while we work.
$rowsToDelete = New-Object System.Collections.ArrayList
$lastRow = $ws.Cells.Item($ws.Rows.Count, 2).End(-4162).Row
for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $text = $cell.Text
    if ($gagesToRemove -contains $text) {
        [void]$rowsToDelete.Add($r)
    }
}

$sorted = $rowsToDelete | Sort-Object -Descending
foreach ($r in $sorted) {
    $ws.Rows.Item($r).Delete() | Out-Null
}

# The old selection (C1:C33) pointed past the new, smaller used range;
# reset it to the default top-left cell now that the sheet only spans
# A1:C30.
$ws.Range("A1").Select() | Out-Null
